$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.674.91"
Set-TextValue "D3" "3.409.18"
Set-TextValue "E3" "  -0.43%  "
Set-TextValue "D5" "569.15"
Set-TextValue "E5" "  -0.73%  "
Set-TextValue "D6" "157.45"
Set-TextValue "E6" "  +0.29%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "D8" "3.412.99"
Set-TextValue "E8" "  -0.44%  "
Set-TextValue "E9" "  -8.76%  "
Set-TextValue "E10" "  +0.94%  "
Set-TextValue "D11" "0.119"
Set-TextValue "E11" "  -4.06%  "
Set-TextValue "E12" "  -4.66%  "
Set-TextValue "D13" "3.995.42"
Set-TextValue "E13" "  -0.53%  "
Set-TextValue "E14" "  +0.12%  "
Set-TextValue "D15" "27.02"
Set-TextValue "E15" "  -3.17%  "
Set-TextValue "E16" "  -8.40%  "
Set-TextValue "D17" "63.757.51"
Set-TextValue "E17" "  -1.37%  "
Set-TextValue "D18" "3.413.38"
Set-TextValue "E18" "  -0.26%  "
Set-TextValue "D19" "6.09"
Set-TextValue "E19" "  -4.51%  "
Set-TextValue "D20" "13.60"
Set-TextValue "E20" "  -2.76%  "
Set-TextValue "D21" "386.25"
Set-TextValue "E21" "  +2.05%  "
Set-TextValue "D22" "7.76"
Set-TextValue "E22" "  -3.57%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  -0.03%  "
Set-TextValue "E24" "  -1.99%  "
Set-TextValue "E25" "  -6.48%  "
Set-TextValue "E26" "  -4.56%  "
Set-TextValue "D27" "9.66"
Set-TextValue "E27" "  -6.15%  "
Set-TextValue "E28" "  +0.13%  "
Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.05%  "
Set-TextValue "D30" "6.09"
Set-TextValue "E30" "  -2.02%  "
Set-TextValue "E31" "  -7.05%  "
Set-TextValue "E32" "  -2.69%  "
Set-TextValue "E33" "  +0.05%  "
Set-TextValue "D34" "22.87"
Set-TextValue "E34" "  -1.21%  "
Set-TextValue "D35" "6.94"
Set-TextValue "E35" "  -4.01%  "
Set-TextValue "E36" "  -6.42%  "
Set-TextValue "D37" "160.72"
Set-TextValue "E37" "  +0.68%  "
Set-TextValue "D38" "0.843"
Set-TextValue "E38" "  +8.83%  "
Set-TextValue "E39" "  -4.63%  "
Set-TextValue "B40" "Hedera"
Set-TextValue "C40" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.0724"
Set-TextValue "E40" "  -5.33%  "
Set-TextValue "B41" "Maker"
Set-TextValue "C41" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D41" "2.785.04"
Set-TextValue "E41" "  -3.08%  "
Set-TextValue "D42" "25.88"
Set-TextValue "E42" "  -3.85%  "
Set-TextValue "E43" "  +0.14%  "
Set-TextValue "D44" "25.91"
Set-TextValue "E44" "  -2.74%  "
Set-TextValue "E45" "  -8.97%  "
Set-TextValue "D46" "4.35"
Set-TextValue "E46" "  -5.87%  "
Set-TextValue "E48" "  +7.77%  "
Set-TextValue "D49" "326.93"
Set-TextValue "E49" "  +1.93%  "
Set-TextValue "E50" "  -5.38%  "
Set-TextValue "E51" "  -5.12%  "
